# Insert a new price record at row 328 of Sheet1, shifting all existing
# rows 328:415 down to 329:416 (matches the upstream diff: dimension
# grows from A1:R415 to A1:R416, and every row from 328 on is pushed
# down by one with the original row 415's data duplicated into the new
# row 416).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 328:415 down to 329:416, leaving a blank row 328 behind.
$ws.Rows(328).Insert()

# Populate the newly-inserted row 328 with the new record.
$ws.Range("A328").Value = 4
$ws.Range("B328").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C328").Value = 'Los Lagos'
$ws.Range("D328").Value = 44943
$ws.Range("E328").Value = 10
$ws.Range("F328").Value = 100112045
$ws.Range("G328").Value = 'Zapallo'
$ws.Range("H328").Value = 'Paine'
$ws.Range("I328").Value = '1a (cosecha)'
$ws.Range("J328").Value = 1200
$ws.Range("K328").Value = 650
$ws.Range("L328").Value = 650
$ws.Range("M328").Value = 650
$ws.Range("N328").Value = '$/kilo (volumen en unidades)'
$ws.Range("O328").Value = "Región de O'Higgins"
$ws.Range("P328").Value = 650
$ws.Range("Q328").Value = 1
$ws.Range("R328").Value = 'Hortaliza'
